$d = $word.ActiveDocument

# Character helpers (avoid relying on literal non-ASCII glyphs in the script
# file so this works regardless of source encoding handling).
$bullet = [char]0x2022
$pm = [char]0x00B1

# Locate the "KEY ACHIEVEMENTS AND IMPACT" heading so we operate on the
# correct bullet list (there is similar-looking text earlier in the
# PROFESSIONAL EXPERIENCE section that must NOT be touched).
$achIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text.Trim()
    if ($t -eq "KEY ACHIEVEMENTS AND IMPACT") {
        $achIdx = $i
        break
    }
}
if ($achIdx -eq -1) {
    throw "Could not find KEY ACHIEVEMENTS AND IMPACT heading"
}

# Bullets start two paragraphs after the heading: heading, "Impact", bullets...
$firstBullet = $achIdx + 2

function Get-ParagraphText($paraIndex) {
    return $d.Paragraphs.Item($paraIndex).Range.Text.Trim()
}

function Set-ParagraphText($paraIndex, $newText) {
    $p = $d.Paragraphs.Item($paraIndex)
    # Replace everything up to (but not including) the paragraph mark so the
    # paragraph's formatting/mark is preserved.
    $r = $d.Range($p.Range.Start, $p.Range.End - 1)
    $r.Text = $newText
}

# Sanity-check the six bullets we expect to find, so we fail loudly instead
# of silently editing the wrong paragraphs if the document layout differs
# from what this script assumes.
$expectedOld1 = $bullet + " Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving classification accuracy from 23% to 64%"
$expectedOld2 = $bullet + " Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%, reducing polling error margins from " + $pm + "4.2% to " + $pm + "2.1%"
$expectedOld3 = $bullet + " Built cloud-based data warehouse solutions on AWS processing billions of records with 99.94% accuracy"
$expectedOld6 = $bullet + " Designed ETL pipelines using PySpark, dbt, and PostgreSQL/PostGIS for large-scale geospatial datasets"

$actual1 = Get-ParagraphText $firstBullet
$actual2 = Get-ParagraphText ($firstBullet + 1)
$actual3 = Get-ParagraphText ($firstBullet + 2)
$actual6 = Get-ParagraphText ($firstBullet + 5)

if ($actual1 -ne $expectedOld1) { throw "Unexpected text at bullet 1: $actual1" }
if ($actual2 -ne $expectedOld2) { throw "Unexpected text at bullet 2: $actual2" }
if ($actual3 -ne $expectedOld3) { throw "Unexpected text at bullet 3: $actual3" }
if ($actual6 -ne $expectedOld6) { throw "Unexpected text at bullet 6: $actual6" }

# Paragraph 1 of the bullet list:
#   "Discovered systematic race coding errors ..." -> "Predictive excellence: ..."
$text1 = $bullet + " Predictive excellence: Achieved 87% voter turnout accuracy vs. 71% industry standard"
Set-ParagraphText $firstBullet $text1

# Paragraph 2:
#   "Achieved 87% prediction accuracy ..." -> "Reduced polling margins from ..."
$idx2 = $firstBullet + 1
$text2 = $bullet + " Reduced polling margins from " + $pm + "4.2% to " + $pm + "2.1%"
Set-ParagraphText $idx2 $text2

# Paragraph 3:
#   "Built cloud-based data warehouse solutions ..." -> "Methodological advancement: ..."
$idx3 = $firstBullet + 2
$text3 = $bullet + " Methodological advancement: Improved segmentation accuracy 34% and survey incidence 28%"
Set-ParagraphText $idx3 $text3

# Paragraphs 4 and 5 ("Built redistricting platform ..." and "Developed
# longitudinal data analysis methods ...") are removed entirely. Delete the
# paragraph at position (firstBullet+3) twice since the list re-indexes after
# each deletion.
$idx4 = $firstBullet + 3
$d.Paragraphs.Item($idx4).Range.Delete()
$d.Paragraphs.Item($idx4).Range.Delete()

# Remaining bullet ("Designed ETL pipelines ...") -> "Reduced polling costs ..."
$text4 = $bullet + " Reduced polling costs while increasing quality"
Set-ParagraphText $idx4 $text4

Write-Host "=== RESULT ==="
for ($i = $achIdx; $i -le $achIdx + 6; $i++) {
    Write-Host "$i => " $d.Paragraphs.Item($i).Range.Text
}
